{"js": "// Finalize Phase 1 Report\n// \"Members: Jonathan, Sriram, Ryan Luna\" -> \"Members: Jonathan Jackson, Sriram Arjula, Ryan Luna\"\n// Add each member's last name right after their first name.\n\nconst body = context.document.body;\n\nasync function insertLastNameAfter(firstName, lastName) {\n  // Idempotent / safe no-op if the last name is already present right after\n  // the first name (e.g. script re-run on an already-edited document).\n  const already = body.search(firstName + \" \" + lastName, { matchCase: true });\n  already.load(\"items\");\n  await context.sync();\n  if (already.items.length > 0) {\n    return false;\n  }\n\n  const found = body.search(firstName, { matchCase: true, matchWholeWord: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    return false; // name not present \u2014 leave document untouched\n  }\n\n  found.items[0].insertText(\" \" + lastName, \"After\");\n  await context.sync();\n  return true;\n}\n\n// Jonathan -> Jonathan Jackson\nawait insertLastNameAfter(\"Jonathan\", \"Jackson\");\n\n// Sriram -> Sriram Arjula\nawait insertLastNameAfter(\"Sriram\", \"Arjula\");\n\n// \"Ryan Luna\" already existed in the document before this edit, so it is left as-is.\n", "ps1": "# Finalize Phase 1 Report\n# \"Members: Jonathan, Sriram, Ryan Luna\" -> \"Members: Jonathan Jackson, Sriram Arjula, Ryan Luna\"\n# Add each member's last name right after their first name.\n\n$d = $word.ActiveDocument\n\nfunction Test-TextExists($needle) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $needle\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.Forward = $true\n    return $rng.Find.Execute()\n}\n\nfunction Insert-LastNameAfter($firstName, $lastName) {\n    # Idempotent / safe no-op if the last name is already present.\n    $already = Test-TextExists ($firstName + \" \" + $lastName)\n    if ($already) {\n        return\n    }\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $firstName\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $true\n    $rng.Find.Forward = $true\n    $found = $rng.Find.Execute()\n    if ($found) {\n        $rng.Collapse(0)  # wdCollapseEnd\n        $rng.InsertAfter(\" \" + $lastName)\n    }\n}\n\n# Jonathan -> Jonathan Jackson\nInsert-LastNameAfter \"Jonathan\" \"Jackson\"\n\n# Sriram -> Sriram Arjula\nInsert-LastNameAfter \"Sriram\" \"Arjula\"\n\n# \"Ryan Luna\" already existed in the document before this edit, so it is left as-is.\n"}
